# "Finished scrapper for all sheets"
# The scraper walked every visible price sheet (US / RMB / HK), refreshing
# the last-seen selection on each tab, and topped up two RMB quotes that
# had not resolved on the previous run (600309 / 600519 switched from the
# Shanghai "SH" suffix used by the old data source to the "SS" suffix used
# by the new one, and now carry real prices instead of blanks).

$wb = $excel.ActiveWorkbook

$hk  = $wb.Worksheets.Item("HK")
$rmb = $wb.Worksheets.Item("RMB")
$us  = $wb.Worksheets.Item("US")

# --- RMB sheet: fill in the two tickers that previously scraped empty ---
$rmb.Range("B8").Value = "600309.SS"
$rmb.Range("C8").Value = 118.9000015258789

$rmb.Range("B9").Value = "600519.SS"
$rmb.Range("C9").Value = 1900

# --- Walk the sheets, leaving each one's last-used cell where the scraper
#     left it. US is touched first, then RMB, finishing on HK which ends
#     up the active tab. ---
$us.Activate()
$us.Range("E15").Select()

$rmb.Activate()
$rmb.Range("D17").Select()

$hk.Activate()
$hk.Range("F11").Select()
